# Apply cryptos-list data refresh (price / 1h-volume columns),
# plus a row swap for EnergySwap <-> RenderToken (rows 45-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to stay a text value (matches the scraped-data
    # inline strings in the source file) instead of letting Excel
    # auto-coerce numeric-looking text (e.g. "529.74") into a number.
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "60.928.23"
Set-TextValue "E2" "  +1.11%  "

Set-TextValue "D3" "2.632.64"
Set-TextValue "E3" "  +1.63%  "

Set-TextValue "E4" "  -0.08%  "

Set-TextValue "D5" "529.74"
Set-TextValue "E5" "  +4.12%  "

Set-TextValue "D6" "155.47"
Set-TextValue "E6" "  +1.12%  "

Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.95%  "

Set-TextValue "D8" "0.589"
Set-TextValue "E8" "  -0.24%  "

Set-TextValue "D9" "6.65"
Set-TextValue "E9" "  -1.03%  "

Set-TextValue "D10" "0.109"
Set-TextValue "E10" "  +5.55%  "

Set-TextValue "D11" "0.351"
Set-TextValue "E11" "  +1.29%  "

Set-TextValue "E12" "  +0.16%  "

Set-TextValue "D13" "3.094.32"
Set-TextValue "E13" "  +1.57%  "

Set-TextValue "D14" "60.920.14"
Set-TextValue "E14" "  +1.12%  "

Set-TextValue "D15" "22.02"
Set-TextValue "E15" "  +2.17%  "

Set-TextValue "D16" "0.0000144"
Set-TextValue "E16" "  +3.36%  "

Set-TextValue "D17" "2.634.49"
Set-TextValue "E17" "  +1.46%  "

Set-TextValue "E18" "  +0.56%  "

Set-TextValue "D19" "353.24"
Set-TextValue "E19" "  +0.19%  "

Set-TextValue "D20" "10.61"
Set-TextValue "E20" "  +0.84%  "

Set-TextValue "E21" "  +2.28%  "

Set-TextValue "D22" "1.00"
Set-TextValue "E22" "  +0.23%  "

Set-TextValue "D23" "61.66"
Set-TextValue "E23" "  +2.14%  "

Set-TextValue "E24" "  +2.28%  "

Set-TextValue "E25" "  +1.45%  "

Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.16%  "

Set-TextValue "D27" "0.0₃0867"
Set-TextValue "E27" "  +3.57%  "

Set-TextValue "E29" "  -0.04%  "

Set-TextValue "D30" "6.16"
Set-TextValue "E30" "  +7.46%  "

Set-TextValue "D31" "19.51"
Set-TextValue "E31" "  +0.65%  "

Set-TextValue "D32" "1.62"
Set-TextValue "E32" "  +4.18%  "

Set-TextValue "D33" "151.09"
Set-TextValue "E33" "  -0.49%  "

Set-TextValue "E34" "  +4.09%  "

Set-TextValue "D35" "1.20"
Set-TextValue "E35" "  +1.84%  "

Set-TextValue "D36" "0.932"
Set-TextValue "E36" "  +11.13%  "

Set-TextValue "D37" "0.886"
Set-TextValue "E37" "  +2.42%  "

Set-TextValue "E38" "  +1.50%  "

Set-TextValue "D39" "3.83"
Set-TextValue "E39" "  +2.16%  "

Set-TextValue "D40" "306.52"
Set-TextValue "E40" "  +3.54%  "

Set-TextValue "D41" "0.641"
Set-TextValue "E41" "  +3.77%  "

Set-TextValue "E42" "  +1.68%  "

Set-TextValue "D43" "0.0562"
Set-TextValue "E43" "  +2.01%  "

Set-TextValue "D44" "0.997"
Set-TextValue "E44" "  -0.05%  "

Set-TextValue "B45" "RenderToken"
Set-TextValue "C45" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D45" "5.03"
Set-TextValue "E45" "  +5.20%  "

Set-TextValue "B46" "EnergySwap"
Set-TextValue "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "19.72"
Set-TextValue "E46" "  +0.05%  "

Set-TextValue "E47" "  +2.74%  "

Set-TextValue "E48" "  +8.04%  "

Set-TextValue "D49" "10.34"
Set-TextValue "E49" "  +0.35%  "

Set-TextValue "D50" "1.982.55"
Set-TextValue "E50" "  -0.40%  "

Set-TextValue "D51" "1.83"
Set-TextValue "E51" "  +2.97%  "
